$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Update status text (shared everywhere "Ready for handoff" is used)
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Add Error Detail messages for row 3 (file 5fddcb5c-f2b6-4949-a993-c27db6d12cf2)
$zhcn.Range("L3").Value = "Handback file name: wwfaog2f.lfh is different with handoff file name: 5fddcb5c-f2b6-4949-a993-c27db6d12cf2.1599b9712162df1f971b5aa296bc5b66c0fbb45a.zh-cn."
$dede.Range("L3").Value = "Handback file name: wwfaog2f.lfh is different with handoff file name: 5fddcb5c-f2b6-4949-a993-c27db6d12cf2.1599b9712162df1f971b5aa296bc5b66c0fbb45a.de-de."
